$d = $word.ActiveDocument

# --- 1. Remove the "Meta description" paragraph that sits right after the
#        title heading (paragraph 2). ---------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$nextPara = $d.Paragraphs.Item(3)
$metaRange = $d.Range($metaPara.Range.Start, $nextPara.Range.Start)
$metaRange.Delete()

# --- 2. Insert a new bold "Play Chilli Heat for Free - Game Review"
#        paragraph right before the closing "Prompt:" paragraph, leaving the
#        paragraph that precedes it (the last bullet item) untouched. -------
$n = $d.Paragraphs.Count
$precedingPara = $d.Paragraphs.Item($n - 1)
$insertionPoint = $precedingPara.Range
$insertionPoint.Collapse(0)

$wordmlNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$precedingXml = "<w:p $wordmlNs><w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Not particularly elaborate graphics</w:t></w:r></w:p>"
$newHeadingXml = "<w:p $wordmlNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Chilli Heat for Free - Game Review</w:t></w:r></w:p>"
$insertionPoint.InsertXML($precedingXml + $newHeadingXml)

# --- 3. Replace the old "Prompt: ..." text of the final paragraph with the
#        new meta-description text, keeping its italic formatting intact. ---
$oldText = "Prompt: Create an appealing feature image for Chilli Heat that is in line with the game's theme and features a happy Maya warrior with glasses. The image should be in a cartoon style. The feature image should include a fun and festive background, perhaps with a Mexican street party, and a cartoon version of the happy Maya warrior as the main focus. The warrior should be holding some chilli peppers and a tequila glass, with a big smile and his signature glasses. In the background, some of the symbols from the game can be included, such as the mariachi, the chihuahua with the tabasco sauce, and the sacks of coins. The overall feel of the image should be colorful and engaging, with a touch of humor to reflect the fun and laid-back nature of the game."
$newText = "Discover the gameplay, betting range, and bonus games of Chilli Heat online slot. Play for free and win up to 1000x your bet."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
